# Update working-hours tracking spreadsheet:
#  - append a new bullet point to the "5.&6.8.2019" entry's notes (C14)
#  - add a new log entry row (15): date 7.8.2019, 4 hours, notes about backend work
#  - grow row heights for the wrapped text (row 14 now taller; new row 15 sized to its text)
#  - move the active-cell selection to C13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the existing note in C14 with an additional bullet line ---
$existingNotes = $ws.Range("C14").Text
$ws.Range("C14").Value = $existingNotes + "`n- Tappelua react dev serverin hot reloadin kanssa. Selvisi, että proxyn käyttäminen rikkoo reloadin jos selaimena Firefox"

# --- add the new log row (row 15) ---
$ws.Range("A15").Value = "7.8.2019"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "Backend: lisätty nodemon, lisätty mongodb tietokanta, tapeltu mongoose+TS interfacet kanssa, jatkettu oauth autorisointia Strava API:iin lisäämällä access tokenien haku kun käyttäjä on antanut luvan tietojen käyttöön ja on saatu Stravan koodi tokeneita varten."

# --- row heights reflecting the longer wrapped text ---
$ws.Rows.Item(14).RowHeight = 102.65
$ws.Rows.Item(15).RowHeight = 77.3

# --- move selection, as the author's cursor ended up on C13 ---
$ws.Range("C13").Select() | Out-Null
